$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a numeric-looking string to be stored as TEXT (matches source
# data, which keeps prices like "0.997" / "3.40" as literal text, not numbers),
# without introducing any new cell style (NumberFormat stays untouched/default).
function Set-TextValue {
    param($addr, $text)
    $r = $ws.Range($addr)
    $escapedForFormula = $text -replace '"', '""'
    $r.Formula = '="' + $escapedForFormula + '"'
    $r.Copy()
    $r.PasteSpecial(-4163)
}

$excel.CutCopyMode = 0

# Row 2
$ws.Range('D2').Value = '67.858.40'
$ws.Range('E2').Value = '  +3.07%  '

# Row 3
$ws.Range('D3').Value = '3.287.35'
$ws.Range('E3').Value = '  -0.44%  '

# Row 4
Set-TextValue 'D4' '0.997'
$ws.Range('E4').Value = '  -0.72%  '

# Row 5
Set-TextValue 'D5' '579.49'
$ws.Range('E5').Value = '  +1.95%  '

# Row 6
Set-TextValue 'D6' '180.81'
$ws.Range('E6').Value = '  -1.70%  '

# Row 7
$ws.Range('E7').Value = '  -0.02%  '

# Row 8
Set-TextValue 'D8' '0.585'
$ws.Range('E8').Value = '  +2.50%  '

# Row 9
$ws.Range('D9').Value = '3.282.63'
$ws.Range('E9').Value = '  -0.38%  '

# Row 10
Set-TextValue 'D10' '0.176'
$ws.Range('E10').Value = '  +0.36%  '

# Row 11
Set-TextValue 'D11' '0.575'
$ws.Range('E11').Value = '  +0.65%  '

# Row 12
Set-TextValue 'D12' '45.72'
$ws.Range('E12').Value = '  -0.86%  '

# Row 13
Set-TextValue 'D13' '0.0000271'
$ws.Range('E13').Value = '  +3.02%  '

# Row 14
Set-TextValue 'D14' '690.77'
$ws.Range('E14').Value = '  +12.90%  '

# Row 15
$ws.Range('D15').Value = '3.812.85'
$ws.Range('E15').Value = '  -0.78%  '

# Row 16
Set-TextValue 'D16' '8.38'
$ws.Range('E16').Value = '  -0.47%  '

# Row 17
$ws.Range('D17').Value = '67.863.54'
$ws.Range('E17').Value = '  +2.93%  '

# Row 19
$ws.Range('D19').Value = '3.290.61'
$ws.Range('E19').Value = '  -0.58%  '

# Row 20
Set-TextValue 'D20' '17.43'
$ws.Range('E20').Value = '  -2.15%  '

# Row 21
Set-TextValue 'D21' '10.84'
$ws.Range('E21').Value = '  -0.79%  '

# Row 22
Set-TextValue 'D22' '0.894'
$ws.Range('E22').Value = '  +0.52%  '

# Row 23
Set-TextValue 'D23' '17.36'
$ws.Range('E23').Value = '  -3.49%  '

# Row 24
Set-TextValue 'D24' '5.18'
$ws.Range('E24').Value = '  +4.48%  '

# Row 25
Set-TextValue 'D25' '97.54'
$ws.Range('E25').Value = '  -3.06%  '

# Row 26
Set-TextValue 'D26' '3.97'
$ws.Range('E26').Value = '  -0.30%  '

# Row 27
Set-TextValue 'D27' '2.75'
$ws.Range('E27').Value = '  +1.74%  '

# Row 28
Set-TextValue 'D28' '5.72'
$ws.Range('E28').Value = '  -3.47%  '

# Row 29
Set-TextValue 'D29' '9.41'
$ws.Range('E29').Value = '  -0.03%  '

# Row 30
Set-TextValue 'D30' '32.71'
$ws.Range('E30').Value = '  +5.78%  '

# Row 31
Set-TextValue 'D31' '8.46'
$ws.Range('E31').Value = '  +0.12%  '

# Row 32
Set-TextValue 'D32' '6.74'
$ws.Range('E32').Value = '  +4.73%  '

# Row 33
Set-TextValue 'D33' '587.72'
$ws.Range('E33').Value = '  +6.61%  '

# Row 34
$ws.Range('D34').Value = '3.901.04'
$ws.Range('E34').Value = '  +2.51%  '

# Row 35
Set-TextValue 'D35' '10.87'
$ws.Range('E35').Value = '  +0.27%  '

# Row 36
$ws.Range('E36').Value = '  +0.95%  '

# Row 39
Set-TextValue 'D39' '55.44'
$ws.Range('E39').Value = '  -0.78%  '

# Row 40
Set-TextValue 'D40' '0.131'
$ws.Range('E40').Value = '  +2.16%  '

# Row 41
Set-TextValue 'D41' '3.26'
$ws.Range('E41').Value = '  +4.06%  '

# Row 42
Set-TextValue 'D42' '2.64'

# Row 43
Set-TextValue 'D43' '32.37'
$ws.Range('E43').Value = '  -0.16%  '

# Row 44
Set-TextValue 'D44' '3.40'
$ws.Range('E44').Value = '  +0.74%  '

# Row 45
$ws.Range('D45').Value = '0.0₃0677'
$ws.Range('E45').Value = '  +0.15%  '

# Row 46
Set-TextValue 'D46' '0.332'
$ws.Range('E46').Value = '  +0.21%  '

# Row 47
Set-TextValue 'D47' '0.0412'
$ws.Range('E47').Value = '  +1.63%  '

# Row 48
$ws.Range('E48').Value = '  +1.77%  '

# Row 49
$ws.Range('E49').Value = '  +0.56%  '

# Row 50
$ws.Range('E50').Value = '  +9.01%  '

# Row 51
Set-TextValue 'D51' '2.51'
$ws.Range('E51').Value = '  +0.57%  '

# Rows 37/38 swapped: Dai and dogwifhat traded rank positions
# Row 37
$ws.Range('B37').Value = 'Dai'
$ws.Range('C37').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D37' '0.996'
$ws.Range('E37').Value = '  -0.12%  '

# Row 38
$ws.Range('B38').Value = 'dogwifhat'
$ws.Range('C38').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 'D38' '3.39'
$ws.Range('E38').Value = '  -9.12%  '

$excel.CutCopyMode = 0
